# CodeSystem-validation-status-cs.xlsx edits
#
# The "Metadata" worksheet holds Property/Value pairs. Three values change:
#   - B7  (Experimental) : was blank -> "false"   (must be literal TEXT, not a
#                            Boolean, so it round-trips as a shared string)
#   - B8  (Date)         : "2025-11-28T14:35:57+00:00" -> "2025-11-30T13:08:37+00:00"
#   - B17 (Description)  : was blank -> "Codes for scientific validation status of metrics"

$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")

# Writing the bare word "false" via .Value/.Value2/.Formula gets auto-typed as
# an Excel Boolean (t="b"), which would change both the stored type and the
# cell's number-format style id. To force genuine text, write it as a string
# formula and then collapse the formula down to its literal text result via
# Copy + PasteSpecial(xlPasteValues) -- this keeps the original cell style
# untouched and yields a plain text cell.
$wsMeta.Range("B7").Formula = '="false"'
$wsMeta.Range("B7").Copy() | Out-Null
$wsMeta.Range("B7").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$wsMeta.Range("B8").Value = "2025-11-30T13:08:37+00:00"

$wsMeta.Range("B17").Value = "Codes for scientific validation status of metrics"
